# Append a new data row (row 58) to Sheet1 with the latest metric reading.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A58").Value = "2025-04-29 07:41:08"
$ws.Range("B58").Value = 173
